$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells stay text (matching source inlineStr cells), avoiding Excel
# auto-converting numeric-looking strings (e.g. "0.4770" -> 0.477).
$textCells = @('D2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'B21', 'C21', 'D21', 'E21', 'B22', 'C22', 'D22', 'E22', 'B23', 'C23', 'D23', 'E23', 'B24', 'C24', 'D24', 'E24', 'B25', 'C25', 'D25', 'E25', 'B26', 'C26', 'D26', 'E26', 'B27', 'C27', 'D27', 'E27', 'B28', 'C28', 'D28', 'E28', 'B29', 'C29', 'D29', 'E29', 'B30', 'C30', 'D30', 'E30', 'B31', 'C31', 'D31', 'E31', 'B32', 'C32', 'D32', 'E32', 'B33', 'C33', 'D33', 'E33', 'B34', 'C34', 'D34', 'E34', 'B35', 'C35', 'D35', 'E35', 'B36', 'C36', 'D36', 'E36', 'B37', 'C37', 'D37', 'E37', 'B38', 'C38', 'D38', 'E38', 'B39', 'C39', 'D39', 'E39', 'B40', 'C40', 'D40', 'E40', 'B41', 'C41', 'D41', 'E41', 'B42', 'C42', 'D42', 'E42', 'B43', 'C43', 'D43', 'E43', 'B44', 'C44', 'D44', 'E44', 'B45', 'C45', 'D45', 'E45', 'B46', 'C46', 'D46', 'E46', 'B47', 'C47', 'D47', 'E47', 'B48', 'C48', 'D48', 'E48', 'B49', 'C49', 'D49', 'E49', 'B50', 'C50', 'D50', 'E50', 'B51', 'C51', 'D51', 'E51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '30.322.64'
$ws.Range('D3').Value = '1.860.90'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '234.08'
$ws.Range('E5').Value = '  -2.17%  '
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '0.4770'
$ws.Range('E7').Value = '  -0.72%  '
$ws.Range('D8').Value = '0.2748'
$ws.Range('E8').Value = '  -3.04%  '
$ws.Range('D9').Value = '0.06442'
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('D10').Value = '1.870.11'
$ws.Range('E10').Value = '  -9.51%  '
$ws.Range('D11').Value = '0.07431'
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').Value = '16.12'
$ws.Range('E12').Value = '  -3.41%  '
$ws.Range('D13').Value = '4.997'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('D14').Value = '86.05'
$ws.Range('E14').Value = '  -2.82%  '
$ws.Range('D15').Value = '0.6327'
$ws.Range('E15').Value = '  -4.99%  '
$ws.Range('D16').Value = '30.304.02'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').Value = '0.9995'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = '231.81'
$ws.Range('E18').Value = '  +3.26%  '
$ws.Range('D19').Value = '12.82'
$ws.Range('E19').Value = '  -4.08%  '
$ws.Range('D20').Value = '0.000007386'
$ws.Range('E20').Value = '  -3.15%  '
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '5.109'
$ws.Range('E22').Value = '  -4.40%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '6.026'
$ws.Range('E23').Value = '  -3.48%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '9.300'
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '167.65'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '17.89'
$ws.Range('E26').Value = '  -3.93%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '1.859'
$ws.Range('E27').Value = '  -5.60%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.1008'
$ws.Range('E28').Value = '  +6.11%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '1.383'
$ws.Range('E29').Value = '  -5.18%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '4.239'
$ws.Range('E30').Value = '  -2.31%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '3.918'
$ws.Range('E31').Value = '  -3.01%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.04904'
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '1.150'
$ws.Range('E33').Value = '  -4.99%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '0.7260'
$ws.Range('E34').Value = '  -3.61%  '
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = '0.9995'
$ws.Range('E35').Value = '  -1.16%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.691'
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.01955'
$ws.Range('E37').Value = '  +6.58%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.630'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '0.9098'
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '1.996'
$ws.Range('E40').Value = '  -4.28%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').Value = '105.80'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '0.9998'
$ws.Range('E42').Value = '  -0.58%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '0.4124'
$ws.Range('E43').Value = '  -3.97%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '5.556'
$ws.Range('E44').Value = '  -5.08%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '7.069'
$ws.Range('E45').Value = '  -5.84%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '61.52'
$ws.Range('E46').Value = '  -6.04%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1209'
$ws.Range('E47').Value = '  -5.90%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '8.825'
$ws.Range('E48').Value = '  -1.33%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.403'
$ws.Range('E49').Value = '  -4.96%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.05614'
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = '33.09'
$ws.Range('E51').Value = '  -2.33%  '
